$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that currently lives on C2 (mailto:juan@example.com).
# The cell keeps its "Hipervínculo" style (s=2) but loses the live link.
$ws.Range("C2").Hyperlinks.Delete()

# Drop the now-unused trailing columns (F:I) entirely so the sheet's
# used range shrinks back down to A1:E2.
$ws.Range("F1:I2").Clear()

# Row 1 - headers
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "locacalizacion"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "id"
$ws.Range("E1").Value = "kind"

# Row 2 - data
$ws.Range("A2").Value = "jorge"
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("E2").Value = 1

# Match the saved selection state (whole populated range selected).
$ws.Range("A1:E2").Select()
